$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row Price (D) and Volume(1h) (E) updates.
# Leading apostrophe forces text entry so values like "1.000" or
# "0.000007770" are not auto-converted to numbers by Excel, matching
# the original inline-string (text) cell type.
$ws.Range("D2").Value = "'29.815.33"
$ws.Range("E2").Value = "'  -1.82%  "
$ws.Range("D3").Value = "'1.889.61"
$ws.Range("E3").Value = "'  -1.78%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'0.7686"
$ws.Range("E5").Value = "'  -5.34%  "
$ws.Range("D6").Value = "'244.52"
$ws.Range("E6").Value = "'  +0.03%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("D8").Value = "'0.3126"
$ws.Range("E8").Value = "'  -4.25%  "
$ws.Range("D9").Value = "'25.24"
$ws.Range("E9").Value = "'  -7.18%  "
$ws.Range("D10").Value = "'0.07202"
$ws.Range("E10").Value = "'  -1.02%  "
$ws.Range("D11").Value = "'0.08101"
$ws.Range("E11").Value = "'  -0.02%  "
$ws.Range("D12").Value = "'0.7658"
$ws.Range("E12").Value = "'  -3.55%  "
$ws.Range("D13").Value = "'5.520"
$ws.Range("E13").Value = "'  +1.74%  "
$ws.Range("D14").Value = "'1.911.16"
$ws.Range("E14").Value = "'  -0.86%  "
$ws.Range("D15").Value = "'92.30"
$ws.Range("E15").Value = "'  -2.45%  "
$ws.Range("D16").Value = "'6.126"
$ws.Range("E16").Value = "'  +0.23%  "
$ws.Range("D17").Value = "'29.834.89"
$ws.Range("E17").Value = "'  -1.82%  "
$ws.Range("E18").Value = "'  -3.22%  "
$ws.Range("D19").Value = "'243.06"
$ws.Range("E19").Value = "'  -3.61%  "
$ws.Range("D20").Value = "'0.000007770"
$ws.Range("E20").Value = "'  -1.37%  "
$ws.Range("D23").Value = "'2.153.08"
$ws.Range("E23").Value = "'  -1.97%  "
$ws.Range("E24").Value = "'  -0.06%  "
$ws.Range("D25").Value = "'0.1560"
$ws.Range("E25").Value = "'  -7.34%  "
$ws.Range("D26").Value = "'9.396"
$ws.Range("E26").Value = "'  -1.61%  "
$ws.Range("D27").Value = "'162.36"
$ws.Range("E27").Value = "'  -3.15%  "
$ws.Range("D28").Value = "'18.74"
$ws.Range("E28").Value = "'  -2.14%  "
$ws.Range("D29").Value = "'2.042"
$ws.Range("E29").Value = "'  -5.70%  "
$ws.Range("D30").Value = "'1.466"
$ws.Range("E30").Value = "'  +6.62%  "
$ws.Range("D31").Value = "'1.550"
$ws.Range("E31").Value = "'  -0.06%  "
$ws.Range("D32").Value = "'4.454"
$ws.Range("E32").Value = "'  +2.26%  "
$ws.Range("D33").Value = "'4.084"
$ws.Range("E33").Value = "'  -1.63%  "
$ws.Range("D34").Value = "'0.05525"
$ws.Range("E34").Value = "'  -2.43%  "
$ws.Range("D35").Value = "'1.257"
$ws.Range("E35").Value = "'  -3.40%  "
$ws.Range("D36").Value = "'0.7478"
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("E37").Value = "'  +0.25%  "
$ws.Range("D38").Value = "'2.632"
$ws.Range("E38").Value = "'  -3.40%  "
$ws.Range("D39").Value = "'0.01921"
$ws.Range("E39").Value = "'  -2.27%  "
$ws.Range("D40").Value = "'2.779"
$ws.Range("E40").Value = "'  -1.58%  "
$ws.Range("D41").Value = "'1.154.80"
$ws.Range("E41").Value = "'  +11.41%  "
$ws.Range("D42").Value = "'73.62"
$ws.Range("E42").Value = "'  -1.15%  "
$ws.Range("D43").Value = "'0.4420"
$ws.Range("E43").Value = "'  -2.07%  "
$ws.Range("D44").Value = "'5.905"
$ws.Range("E44").Value = "'  -1.58%  "
$ws.Range("D45").Value = "'0.8484"
$ws.Range("E45").Value = "'  -1.07%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("D47").Value = "'102.99"
$ws.Range("E47").Value = "'  -0.47%  "
$ws.Range("D48").Value = "'1.884"
$ws.Range("E48").Value = "'  -2.54%  "
$ws.Range("D49").Value = "'9.887"
$ws.Range("E49").Value = "'  -0.77%  "
$ws.Range("D50").Value = "'3.022"
$ws.Range("E50").Value = "'  +0.93%  "
$ws.Range("D51").Value = "'7.445"
$ws.Range("E51").Value = "'  -2.85%  "

# Rows 21/22: Chainlink and Dai swapped ranking positions, with updated data.
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "'  +0.08%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'8.184"
$ws.Range("E22").Value = "'  +1.89%  "
